# Refresh the "cryptos" price table with the latest scrape values.
#
# The sheet keeps Price (D) and Volume(1h) (E) as plain TEXT, not numbers —
# values like "45.415.70" (thousands separators) or "  +0.08%  " (padded
# percent strings) aren't valid Excel numerics anyway, but some refreshed
# prices (e.g. "312.41") *do* look like plain numbers. Left alone, Excel's
# COM layer would auto-coerce those into numeric cells on assignment, which
# would change the stored type. So every D-column cell whose new value
# parses as a number is pre-formatted as Text ("@") before the value is
# written, keeping it a string exactly like its neighbours.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @(
    "D5", "D6", "D7", "D10", "D11", "D12", "D14", "D16",
    "D19", "D22", "D23", "D24", "D27", "D28", "D30", "D31",
    "D32", "D33", "D38", "D39", "D42", "D43", "D44", "D46",
    "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# A1 address -> refreshed value, in sheet (row) order.
$updates = [ordered]@{
    # Bitcoin
    "D2"  = "45.427.38"
    "E2"  = "  +0.07%  "
    # Ethereum
    "D3"  = "2.368.84"
    "E3"  = "  -0.30%  "
    # TetherUSD
    "E4"  = "  -0.10%  "
    # BNB
    "D5"  = "312.41"
    "E5"  = "  -1.67%  "
    # Solana
    "D6"  = "108.17"
    "E6"  = "  -2.95%  "
    # XRP
    "D7"  = "0.629"
    # USDC
    "E8"  = "  -0.03%  "
    # Cardano
    "E9"  = "  -3.29%  "
    # Avalanche
    "D10" = "40.82"
    "E10" = "  -2.60%  "
    # Dogecoin
    "D11" = "0.0917"
    "E11" = "  -1.38%  "
    # Polkadot
    "D12" = "8.47"
    "E12" = "  -2.34%  "
    # TRON
    "E13" = "  +0.81%  "
    # Polygon
    "D14" = "0.978"
    "E14" = "  -3.87%  "
    # WrappedliquidstakedEther2.0
    "D15" = "2.729.06"
    "E15" = "  -0.32%  "
    # Chainlink
    "D16" = "15.32"
    "E16" = "  -2.92%  "
    # WrappedEther
    "D17" = "2.364.79"
    # WrappedBTC
    "D18" = "45.443.29"
    "E18" = "  +0.47%  "
    # InternetComputer(DFINITY)
    "D19" = "14.08"
    "E19" = "  +8.07%  "
    # ShibaInu
    "E20" = "  -1.73%  "
    # Uniswap
    "E21" = "  -5.40%  "
    # Litecoin
    "D22" = "73.26"
    # PancakeSwap
    "D23" = "3.54"
    "E23" = "  -0.15%  "
    # BitcoinCash
    "D24" = "259.54"
    "E24" = "  -3.43%  "
    # ImmutableX
    "E25" = "  +1.78%  "
    # Dai
    "E26" = "  -0.02%  "
    # Cosmos
    "D27" = "11.10"
    "E27" = "  -1.58%  "
    # Filecoin
    "D28" = "7.25"
    "E28" = "  -4.74%  "
    # Toncoin
    "E29" = "  -1.61%  "
    # Hedera
    "D30" = "0.0976"
    "E30" = "  +4.32%  "
    # EthereumClassic
    "D31" = "22.29"
    "E31" = "  -2.85%  "
    # InjectiveProtocol
    "D32" = "36.88"
    "E32" = "  -5.33%  "
    # Monero
    "D33" = "166.30"
    "E33" = "  -2.14%  "
    # WEMIXToken
    "E34" = "  -0.88%  "
    # Kaspa
    "E36" = "  +0.70%  "
    # RenderToken
    "E37" = "  -2.38%  "
    # Row 38 now holds ARBITRUM (was NEARProtocol)
    "B38" = "ARBITRUM"
    "C38" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D38" = "1.90"
    "E38" = "  +9.44%  "
    # Row 39 now holds NEARProtocol (was ARBITRUM)
    "B39" = "NEARProtocol"
    "C39" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D39" = "3.94"
    "E39" = "  +0.43%  "
    # LidoDAOToken
    "E40" = "  -4.82%  "
    # VeChain
    "E41" = "  -3.52%  "
    # BitcoinSV
    "D42" = "98.57"
    "E42" = "  -6.99%  "
    # MultiversX
    "D43" = "69.92"
    "E43" = "  -2.38%  "
    # Algorand
    "D44" = "0.227"
    "E44" = "  -5.51%  "
    # FirstDigitalUSD
    "E45" = "  -0.13%  "
    # Celestia
    "D46" = "12.70"
    "E46" = "  -8.24%  "
    # Maker
    "D47" = "1.811.89"
    "E47" = "  +9.57%  "
    # ordi
    "D48" = "83.64"
    "E48" = "  +5.31%  "
    # THORChain
    "D49" = "5.80"
    "E49" = "  +0.92%  "
    # Row 50 now holds Aave (was FraxShare)
    "B50" = "Aave"
    "C50" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D50" = "111.18"
    "E50" = "  -5.78%  "
    # Row 51 now holds FraxShare (was Aave)
    "B51" = "FraxShare"
    "C51" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D51" = "9.22"
    "E51" = "  +0.80%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
